$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$ws.Range("H17").Value = 479.56924
$ws.Range("J17").Value = 479.56924
$ws.Range("L17").Value = 1438.70772
$ws.Range("N17").Value = -1774.70772

# Row 70 (ALC)
$ws.Range("H70").Value = 2840.8572
$ws.Range("I70").Value = 2600
$ws.Range("J70").Value = 3162
$ws.Range("K70").Value = 7800
$ws.Range("L70").Value = 9486
$ws.Range("M70").Value = -7530
$ws.Range("N70").Value = -10026

# Row 73 (ALC)
$ws.Range("H73").Value = 2840.8572
$ws.Range("I73").Value = 2600
$ws.Range("J73").Value = 3162
$ws.Range("K73").Value = 7800
$ws.Range("L73").Value = 9486
$ws.Range("M73").Value = -6864
$ws.Range("N73").Value = -11358

# Row 86 (ALC)
$ws.Range("H86").Value = 9166.666999999999
$ws.Range("I86").Value = 6000
$ws.Range("J86").Value = 9800
$ws.Range("K86").Value = 6000
$ws.Range("L86").Value = 9800
$ws.Range("M86").Value = -4877
$ws.Range("N86").Value = -12046

# Row 89 (ALC)
$ws.Range("H89").Value = 9166.666999999999
$ws.Range("I89").Value = 6000
$ws.Range("J89").Value = 9800
$ws.Range("K89").Value = 30000
$ws.Range("L89").Value = 49000
$ws.Range("M89").Value = -24384
$ws.Range("N89").Value = -60232

# Row 138 (ALC)
$ws.Range("H138").Value = 513243.53
$ws.Range("I138").Value = 1469.7693
$ws.Range("J138").Value = 738770.9399999999
$ws.Range("K138").Value = 4409.3079
$ws.Range("L138").Value = 2216312.82
$ws.Range("M138").Value = 730.6921000000002
$ws.Range("N138").Value = -2226592.82

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 2814.098
$ws.Range("I32").Value = 2445.694
$ws.Range("J32").Value = 7287.5713
$ws.Range("K32").Value = 2445.694
$ws.Range("L32").Value = 7287.5713
$ws.Range("M32").Value = -2158.694
$ws.Range("N32").Value = -7861.5713

# Row 61 (ARM)
$ws.Range("H61").Value = 1453.0333
$ws.Range("I61").Value = 1284.44
$ws.Range("K61").Value = 1284.44
$ws.Range("M61").Value = -1072.44

# Row 74 (ARM)
$ws.Range("H74").Value = 1576.7059
$ws.Range("I74").Value = 873
$ws.Range("K74").Value = 873
$ws.Range("M74").Value = 1

# Row 77 (ARM)
$ws.Range("H77").Value = 1576.7059
$ws.Range("I77").Value = 873
$ws.Range("K77").Value = 4365
$ws.Range("M77").Value = 3

# Row 136 (ARM)
$ws.Range("H136").Value = 1453.0333
$ws.Range("I136").Value = 1284.44
$ws.Range("K136").Value = 3853.32
$ws.Range("M136").Value = -1303.32

# Row 138 (ARM)
$ws.Range("H138").Value = 56318.332
$ws.Range("J138").Value = 56318.332
$ws.Range("L138").Value = 56318.332
$ws.Range("N138").Value = -66598.33199999999

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (BSM)
$ws.Range("H94").Value = 41667668
$ws.Range("I94").Value = 50000600
$ws.Range("J94").Value = 3010
$ws.Range("K94").Value = 50000600
$ws.Range("L94").Value = 3010
$ws.Range("M94").Value = -50000149
$ws.Range("N94").Value = -3912

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (CRP)
$ws.Range("H22").Value = 588
$ws.Range("I22").Value = 485
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 485
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -135
$ws.Range("N22").Value = -1700

# Row 31 (CRP)
$ws.Range("H31").Value = 1244.5555
$ws.Range("I31").Value = 1215.5
$ws.Range("K31").Value = 1215.5
$ws.Range("M31").Value = -920.5

# Row 32 (CRP)
$ws.Range("H32").Value = 1797.5
$ws.Range("I32").Value = 1600
$ws.Range("J32").Value = 1995
$ws.Range("K32").Value = 1600
$ws.Range("L32").Value = 1995
$ws.Range("M32").Value = -1284
$ws.Range("N32").Value = -2627

# Row 34 (CRP)
$ws.Range("H34").Value = 1244.5555
$ws.Range("I34").Value = 1215.5
$ws.Range("K34").Value = 1215.5
$ws.Range("M34").Value = -1013.5

# Row 35 (CRP)
$ws.Range("H35").Value = 1028.4286
$ws.Range("I35").Value = 949.8333
$ws.Range("J35").Value = 1500
$ws.Range("K35").Value = 949.8333
$ws.Range("L35").Value = 1500
$ws.Range("M35").Value = -655.8333
$ws.Range("N35").Value = -2088

# Row 38 (CRP)
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null

# Row 41 (CRP)
$ws.Range("H41").Value = 11573
$ws.Range("J41").Value = 24832.5
$ws.Range("L41").Value = 24832.5
$ws.Range("N41").Value = -25688.5

# Row 46 (CRP)
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null

$ws = $wb.Worksheets.Item("CUL")
# Row 69 (CUL)
$ws.Range("H69").Value = 2657.25
$ws.Range("I69").Value = 549.5
$ws.Range("J69").Value = 3078.8
$ws.Range("K69").Value = 1648.5
$ws.Range("L69").Value = 9236.400000000001
$ws.Range("M69").Value = -837.5
$ws.Range("N69").Value = -10858.4

# Row 72 (CUL)
$ws.Range("H72").Value = 2657.25
$ws.Range("I72").Value = 549.5
$ws.Range("J72").Value = 3078.8
$ws.Range("K72").Value = 4945.5
$ws.Range("L72").Value = 27709.2
$ws.Range("M72").Value = -889.5
$ws.Range("N72").Value = -35821.2

# Row 131 (CUL)
$ws.Range("H131").Value = 17859970
$ws.Range("J131").Value = 2968
$ws.Range("L131").Value = 8904
$ws.Range("N131").Value = -18984

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 3100
$ws.Range("J80").Value = 6520
$ws.Range("K80").Value = 3100
$ws.Range("L80").Value = 6520
$ws.Range("M80").Value = -2102
$ws.Range("N80").Value = -8516

# Row 83 (GSM)
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 3100
$ws.Range("J83").Value = 6520
$ws.Range("K83").Value = 15500
$ws.Range("L83").Value = 32600
$ws.Range("M83").Value = -10508
$ws.Range("N83").Value = -42584

# Row 102 (GSM)
$ws.Range("H102").Value = 1351.3448
$ws.Range("I102").Value = 1351.3448
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1351.3448
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 270.6551999999999
$ws.Range("N102").Value = $null

# Row 127 (GSM)
$ws.Range("H127").Value = 34210.527
$ws.Range("J127").Value = 34210.527
$ws.Range("L127").Value = 34210.527
$ws.Range("N127").Value = -44130.527

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Range("H40").Value = 3418.2856
$ws.Range("I40").Value = 2152.6155
$ws.Range("J40").Value = 5475
$ws.Range("K40").Value = 2152.6155
$ws.Range("L40").Value = 5475
$ws.Range("M40").Value = -2016.6155
$ws.Range("N40").Value = -5747

# Row 100 (LTW)
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = $null

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (WVR)
$ws.Range("H136").Value = 652.63635
$ws.Range("I136").Value = 448
$ws.Range("J136").Value = 1573.5
$ws.Range("K136").Value = 1344
$ws.Range("L136").Value = 4720.5
$ws.Range("M136").Value = 1206
$ws.Range("N136").Value = -9820.5

# Row 141 (WVR)
$ws.Range("H141").Value = 29724
$ws.Range("J141").Value = 29724
$ws.Range("L141").Value = 29724
$ws.Range("N141").Value = -40084
